# Applies the "Add files via upload" edit:
#  1) Columns G (Columna2/"ext") and H (Columna3/"ext2") were re-entered in
#     thousands (divide raw values by 1000); three H cells (rows 10,11,14)
#     keep an explicit "=n/1000" formula.
#  2) The helper statistics columns PromParo/DesvParo/PromLibros/DesvLibros
#     (old table columns Columna10, Columna11, Columna13, Columna14) were
#     removed entirely, and the two "N" (z-score) columns that used to sit
#     after them (Nparo = old Columna12, Nlibro = old Columna15) were kept
#     but pasted as static values, shifting left to become the new last two
#     sheet columns (O, P). The table itself now only spans B2:N20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Divide the raw "ext"/"ext2" columns (G, H) by 1000.
# ---------------------------------------------------------------------
for ($r = 2; $r -le 20; $r++) {
    $gCell = $ws.Cells.Item($r, 7)   # column G
    $gCell.Value2 = $gCell.Value2 / 1000
}

# H column: most rows get a plain divided value, three rows keep a formula.
$hFormulaRows = @{ 10 = "79104"; 11 = "42812"; 14 = "75356" }
for ($r = 2; $r -le 20; $r++) {
    $hCell = $ws.Cells.Item($r, 8)   # column H
    if ($hFormulaRows.ContainsKey($r)) {
        $hCell.Formula = "=" + $hFormulaRows[$r] + "/1000"
    } else {
        $hCell.Value2 = $hCell.Value2 / 1000
    }
}

# ---------------------------------------------------------------------
# 2) Capture the final "Nparo"/"Nlibro" values (rounded, as they were
#    pasted as values) before the source formula columns disappear.
# ---------------------------------------------------------------------
$nparo = @{
    2 = 1.2977; 3 = -0.7916; 4 = -0.2154; 5 = -1.2278; 6 = 0.78101
    7 = -0.9332; 8 = -0.5614; 9 = 0.30939; 10 = -0.6693; 11 = 0.08082
    12 = 1.10937; 13 = -0.4101; 14 = -0.4713; 15 = 0.24822; 16 = -0.827
    17 = -0.864; 18 = -0.9155; 19 = 2.57574; 20 = 1.48441
}
$nlibro = @{
    2 = 0.82862; 3 = -0.2478; 4 = -0.5294; 5 = -0.5983; 6 = -0.5181
    7 = -0.6728; 8 = 0.3855; 9 = -0.0949; 10 = 3.33176; 11 = 0.58371
    12 = -0.5666; 13 = 0.07803; 14 = 1.53518; 15 = -0.4998; 16 = -0.5863
    17 = 0.02419; 18 = -0.7364; 19 = -0.856; 20 = -0.8606
}

# ---------------------------------------------------------------------
# 3) Shrink the table: drop the six trailing calculated columns
#    (Columna10 PromParo, Columna11 DesvParo, Columna12 Nparo,
#     Columna13 PromLibros, Columna14 DesvLibros, Columna15 Nlibro).
#    ListColumns only delete cleanly from the right-hand end, so remove
#    them highest-index first.
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
for ($i = $lo.ListColumns.Count; $i -ge 14; $i--) {
    $lo.ListColumns.Item($i).Delete()
}

# ---------------------------------------------------------------------
# 4) Re-create the surviving Nparo / Nlibro columns as plain values in
#    O and P (outside the now-smaller table), with their original
#    header labels.
# ---------------------------------------------------------------------
$ws.Range("O1").Value2 = "Nparo"
$ws.Range("P1").Value2 = "Nlibro"

for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 15).Value2 = $nparo[$r]   # column O
    $ws.Cells.Item($r, 16).Value2 = $nlibro[$r]  # column P
}

# Clear any leftover header text beyond column P (old Q1:T1 labels).
$ws.Range("Q1:T1").ClearContents()
